$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 258.75
$ws.Range("I20").Value = 258.75
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 258.75
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -28.75
$ws.Range("N20").ClearContents()

$ws.Range("H28").Value = 608.2727
$ws.Range("I28").Value = 438.44446
$ws.Range("K28").Value = 438.44446
$ws.Range("M28").Value = 46.55554000000001

$ws.Range("H35").Value = 258.75
$ws.Range("I35").Value = 258.75
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 258.75
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = 120.25
$ws.Range("N35").ClearContents()

$ws.Range("H41").Value = 451
$ws.Range("I41").Value = 43.666668
$ws.Range("K41").Value = 43.666668
$ws.Range("M41").Value = 396.333332

$ws.Range("H53").Value = 907.8570999999999
$ws.Range("I53").Value = 1000
$ws.Range("J53").Value = 677.5
$ws.Range("K53").Value = 1000
$ws.Range("L53").Value = 677.5
$ws.Range("M53").Value = -363
$ws.Range("N53").Value = -1951.5

$ws.Range("H55").Value = 234.46153
$ws.Range("I55").Value = 64.833336
$ws.Range("J55").Value = 379.85715
$ws.Range("K55").Value = 64.833336
$ws.Range("L55").Value = 379.85715
$ws.Range("M55").Value = 149.166664
$ws.Range("N55").Value = -807.85715

$ws.Range("H62").Value = 4366.778
$ws.Range("I62").Value = 3869.6667
$ws.Range("J62").Value = 4615.3335
$ws.Range("K62").Value = 3869.6667
$ws.Range("L62").Value = 4615.3335
$ws.Range("M62").Value = -3245.6667
$ws.Range("N62").Value = -5863.3335

$ws.Range("H64").Value = 9541.083000000001
$ws.Range("I64").Value = 6299
$ws.Range("J64").Value = 11856.857
$ws.Range("K64").Value = 6299
$ws.Range("L64").Value = 11856.857
$ws.Range("M64").Value = -6051
$ws.Range("N64").Value = -12352.857

$ws.Range("H65").Value = 4366.778
$ws.Range("I65").Value = 3869.6667
$ws.Range("J65").Value = 4615.3335
$ws.Range("K65").Value = 19348.3335
$ws.Range("L65").Value = 23076.6675
$ws.Range("M65").Value = -16228.3335
$ws.Range("N65").Value = -29316.6675

$ws.Range("H67").Value = 9541.083000000001
$ws.Range("I67").Value = 6299
$ws.Range("J67").Value = 11856.857
$ws.Range("K67").Value = 6299
$ws.Range("L67").Value = 11856.857
$ws.Range("M67").Value = -5441
$ws.Range("N67").Value = -13572.857

$ws.Range("H107").Value = 349.64285
$ws.Range("I107").Value = 325
$ws.Range("J107").Value = 497.5
$ws.Range("K107").Value = 325
$ws.Range("L107").Value = 497.5
$ws.Range("M107").Value = 1595
$ws.Range("N107").Value = -4337.5

$ws.Range("H113").Value = 3090
$ws.Range("J113").Value = 2816.6667
$ws.Range("L113").Value = 2816.6667
$ws.Range("N113").Value = -9324.6667

$ws.Range("H116").Value = 6778.857
$ws.Range("J116").Value = 7074.875
$ws.Range("L116").Value = 7074.875
$ws.Range("N116").Value = -13958.875

$ws.Range("H138").Value = 6040.647
$ws.Range("I138").Value = 5928
$ws.Range("J138").Value = 6119.5
$ws.Range("K138").Value = 17784
$ws.Range("L138").Value = 18358.5
$ws.Range("M138").Value = -12644
$ws.Range("N138").Value = -28638.5

$ws.Range("H141").Value = 1999.5
$ws.Range("I141").Value = 1999.5
$ws.Range("K141").Value = 5998.5
$ws.Range("M141").Value = -818.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3106.2856
$ws.Range("I2").Value = 177
$ws.Range("K2").Value = 177
$ws.Range("M2").Value = -64

$ws.Range("H45").Value = 3749.5
$ws.Range("I45").Value = 1499
$ws.Range("J45").Value = 6000
$ws.Range("K45").Value = 1499
$ws.Range("L45").Value = 6000
$ws.Range("M45").Value = -1122
$ws.Range("N45").Value = -6754

$ws.Range("H61").Value = 12998.25
$ws.Range("I61").Value = 14997.667
$ws.Range("K61").Value = 14997.667
$ws.Range("M61").Value = -14785.667

$ws.Range("H63").Value = 2278.2222
$ws.Range("I63").Value = 1500.625
$ws.Range("K63").Value = 1500.625
$ws.Range("M63").Value = -814.625

$ws.Range("H66").Value = 2278.2222
$ws.Range("I66").Value = 1500.625
$ws.Range("K66").Value = 7503.125
$ws.Range("M66").Value = -4071.125

$ws.Range("H74").Value = 1470.3572
$ws.Range("I74").Value = 1414.5454
$ws.Range("K74").Value = 1414.5454
$ws.Range("M74").Value = -540.5454

$ws.Range("H77").Value = 1470.3572
$ws.Range("I77").Value = 1414.5454
$ws.Range("K77").Value = 7072.727
$ws.Range("M77").Value = -2704.727

$ws.Range("H116").Value = 3106.2856
$ws.Range("I116").Value = 177
$ws.Range("K116").Value = 177
$ws.Range("M116").Value = 2117

$ws.Range("H122").Value = 2504.75
$ws.Range("I122").Value = 2504.75
$ws.Range("K122").Value = 7514.25
$ws.Range("M122").Value = -5064.25

$ws.Range("H132").Value = 3548.5
$ws.Range("I132").Value = 3548.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10645.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8115.5
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 12998.25
$ws.Range("I136").Value = 14997.667
$ws.Range("K136").Value = 44993.001
$ws.Range("M136").Value = -42443.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3106.2856
$ws.Range("I3").Value = 177
$ws.Range("K3").Value = 177
$ws.Range("M3").Value = -63

$ws.Range("H64").Value = 156.28572
$ws.Range("I64").Value = 98.666664
$ws.Range("K64").Value = 98.666664
$ws.Range("M64").Value = 126.333336

$ws.Range("H67").Value = 156.28572
$ws.Range("I67").Value = 98.666664
$ws.Range("K67").Value = 98.666664
$ws.Range("M67").Value = 681.333336

$ws.Range("H134").Value = 4997
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2728.7778
$ws.Range("I31").Value = 2728.7778
$ws.Range("K31").Value = 2728.7778
$ws.Range("M31").Value = -2433.7778

$ws.Range("H34").Value = 2728.7778
$ws.Range("I34").Value = 2728.7778
$ws.Range("K34").Value = 2728.7778
$ws.Range("M34").Value = -2526.7778

$ws.Range("H35").Value = 3647.3333
$ws.Range("I35").Value = 3647.3333
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 3647.3333
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -3353.3333
$ws.Range("N35").ClearContents()

$ws.Range("H94").Value = 2884.1667
$ws.Range("I94").Value = 2457.1428
$ws.Range("J94").Value = 3482
$ws.Range("K94").Value = 2457.1428
$ws.Range("L94").Value = 3482
$ws.Range("M94").Value = -2006.1428
$ws.Range("N94").Value = -4384

$ws.Range("H105").Value = 1953.909
$ws.Range("I105").Value = 1398.8
$ws.Range("K105").Value = 1398.8
$ws.Range("M105").Value = 348.2

$ws.Range("H132").Value = 6653.8887
$ws.Range("I132").Value = 4298.7144
$ws.Range("K132").Value = 12896.1432
$ws.Range("M132").Value = -10366.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2541.4285
$ws.Range("I132").Value = 1299.3334
$ws.Range("J132").Value = 2880.182
$ws.Range("K132").Value = 11694.0006
$ws.Range("L132").Value = 25921.638
$ws.Range("M132").Value = -9164.000599999999
$ws.Range("N132").Value = -30981.638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1993.75
$ws.Range("I22").Value = 1500
$ws.Range("K22").Value = 1500
$ws.Range("M22").Value = -1205

$ws.Range("H27").Value = 1993.75
$ws.Range("I27").Value = 1500
$ws.Range("K27").Value = 1500
$ws.Range("M27").Value = -1393

$ws.Range("H32").Value = 1006
$ws.Range("I32").Value = 1006
$ws.Range("K32").Value = 1006
$ws.Range("M32").Value = -689

$ws.Range("H55").Value = 861.53845
$ws.Range("I55").Value = 100
$ws.Range("J55").Value = 925
$ws.Range("K55").Value = 100
$ws.Range("L55").Value = 925
$ws.Range("M55").Value = 73
$ws.Range("N55").Value = -1271

$ws.Range("H61").Value = 4500
$ws.Range("I61").Value = 4500
$ws.Range("K61").Value = 4500
$ws.Range("M61").Value = -4298

$ws.Range("H100").Value = 1854
$ws.Range("I100").Value = 1496.3334
$ws.Range("K100").Value = 1496.3334
$ws.Range("M100").Value = -955.3334

$ws.Range("H113").Value = 4500
$ws.Range("I113").Value = 4500
$ws.Range("K113").Value = 4500
$ws.Range("M113").Value = -2330

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 772.5
$ws.Range("I113").Value = 772.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2317.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -147.5
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 1915.8889
$ws.Range("I132").Value = 1915.8889
$ws.Range("K132").Value = 5747.6667
$ws.Range("M132").Value = -3217.6667
